$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows for the MuSCs sending-cluster group (rows 8-10 in the old layout),
# which dropped out of this TPM run; dimension becomes A1:T7.
$ws.Rows("8:10").Delete() | Out-Null

# Updated values (TPM recompute) for the remaining rows 2-7, column by column (A..T).
$rows = @(
    @{A="ECs"; B="Spon1"; C="Lrp8"; D="FAPs"; E=3; F=1; G=0.392285; H=1.176855; I=0.0216687106565248; J=0.0216687106565248; K=3; L=1; M=0.1124773333333333; N=0.337432; O=0.7871437602495106; P=0.7871437602495107; Q=0.04412317070666666; R=0.3971085363599999; S=0.01705639038593558; T=0.01705639038593558}
    @{A="ECs"; B="Spon1"; C="Lrp8"; D="MuSCs"; E=3; F=1; G=0.392285; H=1.176855; I=0.0216687106565248; J=0.0216687106565248; K=2; L=0.6666666666666666; M=0.03041566666666666; N=0.09124699999999999; O=0.2128562397504893; P=0.2128562397504893; Q=0.01193160979833333; R=0.107384488185; S=0.004612320270589226; T=0.004612320270589226}
    @{A="FAPs"; B="Spon1"; C="Lrp8"; D="FAPs"; E=3; F=1; G=13.91986866666667; H=41.75960600000001; I=0.7688940604785444; J=0.7688940604785441; K=3; L=1; M=0.1124773333333333; N=0.337432; O=0.7871437602495106; P=0.7871437602495107; Q=1.565669707976889; R=14.091027371792; S=0.605230161998596; T=0.6052301619985959}
    @{A="FAPs"; B="Spon1"; C="Lrp8"; D="MuSCs"; E=3; F=1; G=13.91986866666667; H=41.75960600000001; I=0.7688940604785444; J=0.7688940604785441; K=2; L=0.6666666666666666; M=0.03041566666666666; N=0.09124699999999999; O=0.2128562397504893; P=0.2128562397504893; Q=0.4233820854091111; R=3.810438768682; S=0.1636638984799482; T=0.1636638984799482}
    @{A="MuSCs"; B="Spon1"; C="Lrp8"; D="FAPs"; E=3; F=1; G=3.7916; H=11.3748; I=0.209437228864931; J=0.209437228864931; K=3; L=1; M=0.1124773333333333; N=0.337432; O=0.7871437602495106; P=0.7871437602495107; Q=0.4264690570666667; R=3.8382215136; S=0.1648572078649791; T=0.1648572078649791}
    @{A="MuSCs"; B="Spon1"; C="Lrp8"; D="MuSCs"; E=3; F=1; G=3.7916; H=11.3748; I=0.209437228864931; J=0.209437228864931; K=2; L=0.6666666666666666; M=0.03041566666666666; N=0.09124699999999999; O=0.2128562397504893; P=0.2128562397504893; Q=0.1153240417333333; R=1.0379163756; S=0.04458002099995185; T=0.04458002099995185}
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$r = 2
foreach ($row in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $row[$col]
    }
    $r++
}